# Adds a new data row (row 4) to the "Artfynd" sheet, mirroring the shape
# of the existing rows (2 and 3), and expands the sheet's used range from
# A1:AY3 to A1:AY4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- plain numeric cells -------------------------------------------------
$ws.Range("A4").Value = 112079249
$ws.Range("B4").Value = 89646
$ws.Range("E4").Value = 65
$ws.Range("Q4").Value = 613881.387574179
$ws.Range("R4").Value = 7034405.898391382
$ws.Range("S4").Value = 20

# --- plain text cells ------------------------------------------------------
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "VU"
$ws.Range("F4").Value = "Fläckporing"
$ws.Range("G4").Value = "Anthoporia albobrunnea"
$ws.Range("H4").Value = "(Romell) Karasiński & Niemelä"
$ws.Range("P4").Value = "Älggårdshöjden, Ång"
$ws.Range("T4").Value = "Västernorrland"
$ws.Range("U4").Value = "Sollefteå"
$ws.Range("V4").Value = "Ångermanland"
$ws.Range("W4").Value = "Resele"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AW4").Value = "Daniel Rutschman"
$ws.Range("AX4").Value = "Daniel Rutschman"

# Date-looking text must be forced to text (leading apostrophe) so Excel
# doesn't coerce it into a date serial number.
$ws.Range("Y4").Value = "'2023-09-13"
$ws.Range("Y4").Style = "Normal"
$ws.Range("AA4").Value = "'2023-09-13"
$ws.Range("AA4").Style = "Normal"

# --- boolean cells -----------------------------------------------------
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false

# --- empty (but present) text cells -------------------------------------
# A leading apostrophe forces these to be stored as empty text cells
# (matching t="inlineStr" with no content in the source) instead of being
# left as untouched/blank cells; resetting the style afterwards drops the
# "quote prefix" formatting flag that entering a leading apostrophe adds.
$ws.Range("I4").Value = "'"
$ws.Range("I4").Style = "Normal"
$ws.Range("AT4").Value = "'"
$ws.Range("AT4").Style = "Normal"
$ws.Range("AY4").Value = "'"
$ws.Range("AY4").Style = "Normal"
